# cwl_sources.xlsx - "feat(CWL): complete csharp eval & compiler support"
#
# Adds 12 new localization rows (id / text_JP / text) describing the new
# Roslyn-based C# script compiler + a drama `call` failure message, plus a
# block of blank spacer rows, matching the author's spreadsheet edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-RowValues {
    param(
        [int]$Row,
        [string]$A,
        [string]$C,
        [string]$D
    )
    if ($A -ne $null) { $ws.Range("A$Row").Value2 = $A }
    if ($C -ne $null) { $ws.Range("C$Row").Value2 = $C }
    if ($D -ne $null) { $ws.Range("D$Row").Value2 = $D }
}

function Copy-RowFormat {
    param(
        [int]$FromRow,
        [int]$ToRow
    )
    $ws.Range("A$FromRow`:D$FromRow").Copy()
    $ws.Range("A$ToRow`:D$ToRow").PasteSpecial($xlPasteFormats)
}

function Copy-CellFormat {
    param(
        [string]$FromCell,
        [string]$ToCell
    )
    $ws.Range($FromCell).Copy()
    $ws.Range($ToCell).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Grab "donor" cells for styles that are about to be overwritten further
# down (129/132 currently carry style 25 on C/D - that style is reused by
# several of the new rows, so snapshot it onto a scratch area first).
# ---------------------------------------------------------------------
Copy-CellFormat "C129" "C200"
Copy-CellFormat "A6"   "A201"

# ---------------------------------------------------------------------
# Re-style existing rows 129, 130, 132 (content unchanged, only format)
# ---------------------------------------------------------------------
$ws.Rows.Item(129).RowHeight = 23.25
Copy-RowFormat 131 129

$ws.Rows.Item(130).RowHeight = 23.25
Copy-RowFormat 131 130

$ws.Rows.Item(132).RowHeight = 23.25
Copy-CellFormat "A131" "C132"
Copy-CellFormat "A131" "D132"

# ---------------------------------------------------------------------
# Fill in the values for the already-blank (style 24) rows 133 / 134
# ---------------------------------------------------------------------
Set-RowValues 133 "cwl_log_csc_roslyn" "Roslyn コンパイラを使用しています {0}" "Roslyn 编译器 {0}"
Set-RowValues 134 "cwl_log_csc_package" "{1} からパッケージ {0} をコンパイルしています" "正在编译包 {0} << {1}"

# ---------------------------------------------------------------------
# Brand-new rows 135-140 (the feature content)
# ---------------------------------------------------------------------
Set-RowValues 135 "cwl_log_csc_eval" "スクリプトをコンパイルしています`n{0}" "正在编译脚本`n{0}"
$ws.Rows.Item(135).RowHeight = 46.5
Copy-RowFormat 131 135
Copy-CellFormat "C200" "C135"
Copy-CellFormat "C200" "D135"

Set-RowValues 136 "cwl_error_cs_disabled" "スクリプトコンパイラは無効化されています" "脚本编译器已被禁用"
$ws.Rows.Item(136).RowHeight = 23.25
Copy-RowFormat 131 136
Copy-CellFormat "A201" "A136"

Set-RowValues 137 "cwl_log_csc_scripts" "{1} 個のスクリプトファイルをコンパイルしています '{0}'" "正在编译 {1} 个脚本文件 '{0}'"
$ws.Rows.Item(137).RowHeight = 23.25
Copy-RowFormat 131 137

Set-RowValues 138 "cwl_error_csc_diag" "「{0}」のコンパイルに失敗しました:`n{1}" "编译 '{0}' 失败:`n{1}"
$ws.Rows.Item(138).RowHeight = 46.5
Copy-RowFormat 131 138
Copy-CellFormat "A201" "A138"
Copy-CellFormat "C200" "C138"
Copy-CellFormat "C200" "D138"

Set-RowValues 139 "cwl_error_cs_frozen" "スクリプト状態「{0}」は凍結されています" "脚本状态 '{0}' 已被冻结"
$ws.Rows.Item(139).RowHeight = 23.25
Copy-RowFormat 131 139
Copy-CellFormat "A201" "A139"

Set-RowValues 140 "cwl_warn_drama_call_ex" "call failure: '{0}'`n{1}" "调用失败: '{0}'`n{1}"
$ws.Rows.Item(140).RowHeight = 46.5
Copy-RowFormat 131 140
Copy-CellFormat "C200" "C140"
Copy-CellFormat "C200" "D140"
# A140 gets its own brand-new font (Cascadia Code, themed accent color) -
# closest reproduction available through the exposed object model.
$ws.Range("A140").Font.Name = "Cascadia Code"
$ws.Range("A140").Font.Size = 16
$ws.Range("A140").Font.ThemeColor = 8
$ws.Range("A140").Font.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# Trailing blank spacer rows 141-153 (style 24, ht 23.25, all 4 columns)
# ---------------------------------------------------------------------
for ($r = 141; $r -le 153; $r++) {
    Copy-RowFormat 131 $r
    $ws.Rows.Item($r).RowHeight = 23.25
}

# ---------------------------------------------------------------------
# Clean up scratch donor cells used above
# ---------------------------------------------------------------------
$ws.Range("A200:D201").Clear()

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# Column layout: id column narrows slightly, filter column is now hidden
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 51.45
$ws.Columns.Item(2).Hidden = $true

# ---------------------------------------------------------------------
# Sheet view: scroll down to the newly-added content and select D140
# ---------------------------------------------------------------------
$ws.Range("D140").Select()
